$wb = $excel.ActiveWorkbook

# The F-column ("想去人数" / "wish-to-go count") values were updated in both
# the "展览" and "全部类型" worksheets (they contain duplicated data).
$sheetNames = @("展览", "全部类型")

# Map of row -> new F value
$updates = @{
    2  = 6518
    5  = 44
    6  = 1953
    7  = 1492
    9  = 998
    10 = 361
    11 = 4
    12 = 5621
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
